# Add new column O to the "Suivi" sheet, mirroring column N's formatting,
# with O1 containing a new timestamp header and O2:O100 copied from N2:N100
# (same price values), and O101:O204 left blank like N101:N204.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column N (rows 1-204) into column O so formatting/types match exactly.
$ws.Range("N1:N204").Copy($ws.Range("O1:O204"))

# Header cell gets the new timestamp for this snapshot.
$ws.Range("O1").Value2 = "2026-01-28 06:20:14"
